$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.193646430969238
$ws.Range("B1").Value = 2.588683843612671
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 2.192995071411133
$ws.Range("E1").Value = 1.179575204849243
